{"js": "// Update the answer cells in the \"two-digit divided by one-digit\" table.\n// The document holds a single 5-column table; only every 4th row (0, 4, 8,\n// 12, 16 \u2014 i.e. rows 1, 5, 9, 13, 17 of the 20-row grid) carries the\n// \"xx\u00f7y=zz, r\" answers; the rows in between are blank spacer rows.\n// We address cells positionally (row, col) rather than by text search,\n// because the source text \"58\u00f78=7, 2\" is not unique (it occurs twice, and\n// each occurrence must become a different value).\n\nconst table = context.document.body.tables.getFirst();\n\n// Each entry: [rowIndex, colIndex, newText]\nconst updates = [\n  [0, 0, \"48\u00f74=12, 0\"],\n  [0, 1, \"14\u00f73=4, 2\"],\n  [0, 2, \"90\u00f79=10, 0\"],\n  [0, 3, \"27\u00f72=13, 1\"],\n  [0, 4, \"35\u00f72=17, 1\"],\n\n  [4, 0, \"79\u00f72=39, 1\"],\n  [4, 1, \"93\u00f78=11, 5\"],\n  [4, 2, \"87\u00f72=43, 1\"],\n  [4, 3, \"17\u00f74=4, 1\"],\n  [4, 4, \"66\u00f72=33, 0\"],\n\n  [8, 0, \"69\u00f72=34, 1\"],\n  [8, 1, \"65\u00f73=21, 2\"],\n  [8, 2, \"83\u00f75=16, 3\"],\n  [8, 3, \"78\u00f72=39, 0\"],\n  [8, 4, \"14\u00f74=3, 2\"],\n\n  [12, 0, \"37\u00f79=4, 1\"],\n  [12, 1, \"17\u00f75=3, 2\"],\n  [12, 2, \"59\u00f78=7, 3\"],\n  [12, 3, \"73\u00f74=18, 1\"],\n  [12, 4, \"65\u00f76=10, 5\"],\n\n  [16, 0, \"33\u00f72=16, 1\"],\n  [16, 1, \"22\u00f78=2, 6\"],\n  [16, 2, \"82\u00f74=20, 2\"],\n  [16, 3, \"98\u00f78=12, 2\"],\n  [16, 4, \"24\u00f78=3, 0\"],\n];\n\nfor (const [row, col, text] of updates) {\n  table.getCell(row, col).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the answer cells in the \"two-digit divided by one-digit\" table.\n# The document holds a single 5-column table; only every 4th row (Word COM\n# rows 1, 5, 9, 13, 17 of the 20-row grid) carries the \"xx\u00f7y=zz, r\" answers,\n# the rows in between are blank spacer rows.\n# We address cells positionally (row, column) rather than by Find/Replace,\n# because the source text \"58\u00f78=7, 2\" is not unique (it occurs twice, and\n# each occurrence must become a different value).\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$updates = @(\n  @{r=1;  c=1; v=\"48\u00f74=12, 0\"},\n  @{r=1;  c=2; v=\"14\u00f73=4, 2\"},\n  @{r=1;  c=3; v=\"90\u00f79=10, 0\"},\n  @{r=1;  c=4; v=\"27\u00f72=13, 1\"},\n  @{r=1;  c=5; v=\"35\u00f72=17, 1\"},\n\n  @{r=5;  c=1; v=\"79\u00f72=39, 1\"},\n  @{r=5;  c=2; v=\"93\u00f78=11, 5\"},\n  @{r=5;  c=3; v=\"87\u00f72=43, 1\"},\n  @{r=5;  c=4; v=\"17\u00f74=4, 1\"},\n  @{r=5;  c=5; v=\"66\u00f72=33, 0\"},\n\n  @{r=9;  c=1; v=\"69\u00f72=34, 1\"},\n  @{r=9;  c=2; v=\"65\u00f73=21, 2\"},\n  @{r=9;  c=3; v=\"83\u00f75=16, 3\"},\n  @{r=9;  c=4; v=\"78\u00f72=39, 0\"},\n  @{r=9;  c=5; v=\"14\u00f74=3, 2\"},\n\n  @{r=13; c=1; v=\"37\u00f79=4, 1\"},\n  @{r=13; c=2; v=\"17\u00f75=3, 2\"},\n  @{r=13; c=3; v=\"59\u00f78=7, 3\"},\n  @{r=13; c=4; v=\"73\u00f74=18, 1\"},\n  @{r=13; c=5; v=\"65\u00f76=10, 5\"},\n\n  @{r=17; c=1; v=\"33\u00f72=16, 1\"},\n  @{r=17; c=2; v=\"22\u00f78=2, 6\"},\n  @{r=17; c=3; v=\"82\u00f74=20, 2\"},\n  @{r=17; c=4; v=\"98\u00f78=12, 2\"},\n  @{r=17; c=5; v=\"24\u00f78=3, 0\"}\n)\n\nforeach ($u in $updates) {\n  $t.Cell($u.r, $u.c).Range.Text = $u.v\n}\n"}
